$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C5").Value = "Robot3"
$ws.Range("D5").Value = 8.67
$ws.Range("E5").Value = 13
$ws.Range("D6").Value = 8
$ws.Range("E6").Value = 12
